$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# --- Remove two stray explanation rows -------------------------------------
# Row 77 "APPARTENENZA" (S-row) in the "9) tipologia di lavoro" table is a
# duplicate/extraneous line -> delete entire row, shifting everything below
# up by one.
$ws.Rows.Item(77).Delete()

# Row 83 "SPECIFICAZIONE" (S-row, originally row 84 before the previous
# delete shifted it up) in the "10) materiale" table is likewise extraneous
# -> delete entire row.
$ws.Rows.Item(83).Delete()

# --- Fix the two "Totale" summary strings (typo: ...002 -> ...001) ---------
# Both cells hold rich text: a bold "Totale" run followed by a plain run
# with the rest of the sentence. Re-assign the value and re-apply the rich
# text formatting so the bold prefix is preserved. (B84 is updated first so
# the rebuilt shared-string table keeps the same index order as before.)
$r2 = $ws.Range("B84")
$r2.Value = "Totale: 1.000.001L x 1 all'anno = 1.000.001L all'anno = 1.000.001L all'anno"
$r2.Characters(1, 6).Font.Bold = $true
$r2.Characters(7, 69).Font.Bold = $false

$r1 = $ws.Range("B78")
$r1.Value = "Totale: 50.001L x 1 all'anno = 50.001L all'anno = 50.001 all'anno"
$r1.Characters(1, 6).Font.Bold = $true
$r1.Characters(7, 59).Font.Bold = $false

# --- Update print area / selection to the region now being reviewed --------
$ws.PageSetup.PrintArea = 'B81:E84'
$ws.Range("B81:E84").Select()
